# Refresh the cryptos price/volume snapshot (GitHub Actions scrape update).
# For each touched cell we set the literal text exactly as scraped. Cells whose
# new text is a plain decimal number (e.g. "498.40") are written with a leading
# "'" so Excel keeps them as text (matching the original inlineStr cells) instead
# of silently re-parsing them as numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '56.850.44'
$ws.Range('E2').Value = '  -1.75%  '
# Row 3
$ws.Range('D3').Value = '2.980.05'
$ws.Range('E3').Value = '  -2.18%  '
# Row 4
$ws.Range('E4').Value = '  +0.00%  '
# Row 5
$ws.Range('D5').Value = '''498.40'
$ws.Range('E5').Value = '  -5.10%  '
# Row 6
$ws.Range('E6').Value = '  -3.33%  '
# Row 7
$ws.Range('E7').Value = '  +0.09%  '
# Row 8
$ws.Range('D8').Value = '''0.429'
$ws.Range('E8').Value = '  -3.48%  '
# Row 9
$ws.Range('E9').Value = '  -4.22%  '
# Row 10
$ws.Range('E10').Value = '  -4.75%  '
# Row 11
$ws.Range('D11').Value = '''0.357'
$ws.Range('E11').Value = '  -3.36%  '
# Row 12
$ws.Range('D12').Value = '3.495.85'
$ws.Range('E12').Value = '  -2.07%  '
# Row 13
$ws.Range('E13').Value = '  -2.33%  '
# Row 14
$ws.Range('D14').Value = '''26.02'
$ws.Range('E14').Value = '  -2.82%  '
# Row 15
$ws.Range('D15').Value = '''0.0000159'
$ws.Range('E15').Value = '  -7.34%  '
# Row 16
$ws.Range('D16').Value = '56.971.91'
$ws.Range('E16').Value = '  -1.48%  '
# Row 17
$ws.Range('D17').Value = '''6.05'
$ws.Range('E17').Value = '  -2.55%  '
# Row 18
$ws.Range('D18').Value = '2.988.40'
$ws.Range('E18').Value = '  -1.92%  '
# Row 19
$ws.Range('D19').Value = '''12.58'
$ws.Range('E19').Value = '  -2.66%  '
# Row 20
$ws.Range('D20').Value = '''7.85'
$ws.Range('E20').Value = '  -4.27%  '
# Row 21
$ws.Range('D21').Value = '''319.50'
$ws.Range('E21').Value = '  -6.54%  '
# Row 22
$ws.Range('E22').Value = '  -0.15%  '
# Row 23
$ws.Range('E23').Value = '  -0.18%  '
# Row 24
$ws.Range('E24').Value = '  -2.13%  '
# Row 25
$ws.Range('D25').Value = '''63.52'
$ws.Range('E25').Value = '  -2.73%  '
# Row 26
$ws.Range('E26').Value = '  -0.06%  '
# Row 27
$ws.Range('D27').Value = '''0.163'
$ws.Range('E27').Value = '  -5.70%  '
# Row 28
$ws.Range('D28').Value = '0.0₃0888'
$ws.Range('E28').Value = '  -8.25%  '
# Row 29
$ws.Range('D29').Value = '''6.51'
$ws.Range('E29').Value = '  -6.89%  '
# Row 30
$ws.Range('D30').Value = '''7.05'
$ws.Range('E30').Value = '  -3.38%  '
# Row 31
$ws.Range('E31').Value = '  -5.33%  '
# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '''1.16'
$ws.Range('E32').Value = '  -6.68%  '
# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''20.18'
$ws.Range('E33').Value = '  -4.24%  '
# Row 34
$ws.Range('D34').Value = '''155.34'
$ws.Range('E34').Value = '  -0.75%  '
# Row 35
$ws.Range('D35').Value = '''4.55'
$ws.Range('E35').Value = '  -3.74%  '
# Row 36
$ws.Range('E36').Value = '  -2.36%  '
# Row 37
$ws.Range('E37').Value = '  -7.13%  '
# Row 38
$ws.Range('D38').Value = '''24.24'
$ws.Range('E38').Value = '  -6.13%  '
# Row 39
$ws.Range('D39').Value = '''0.0663'
$ws.Range('E39').Value = '  -4.62%  '
# Row 40
$ws.Range('D40').Value = '3.012.17'
$ws.Range('E40').Value = '  -2.22%  '
# Row 41
$ws.Range('D41').Value = '''37.68'
$ws.Range('E41').Value = '  -0.07%  '
# Row 42
$ws.Range('E42').Value = '  -0.08%  '
# Row 43
$ws.Range('D43').Value = '''3.73'
$ws.Range('E43').Value = '  -3.14%  '
# Row 44
$ws.Range('D44').Value = '''0.641'
$ws.Range('E44').Value = '  -3.16%  '
# Row 45
$ws.Range('D45').Value = '2.201.58'
$ws.Range('E45').Value = '  -5.31%  '
# Row 46
$ws.Range('E46').Value = '  -6.46%  '
# Row 47
$ws.Range('D47').Value = '''5.94'
$ws.Range('E47').Value = '  -1.72%  '
# Row 48
$ws.Range('D48').Value = '''0.936'
$ws.Range('E48').Value = '  -9.25%  '
# Row 49
$ws.Range('E49').Value = '  -5.75%  '
# Row 50
$ws.Range('D50').Value = '''19.22'
$ws.Range('E50').Value = '  -4.39%  '
# Row 51
$ws.Range('D51').Value = '''1.79'
$ws.Range('E51').Value = '  -12.01%  '
